$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("C2").Copy($ws.Range("C3"))
$ws.Range("D2").Copy($ws.Range("D3"))
$ws.Range("E2").Copy($ws.Range("E3"))
$ws.Range("G2").Copy($ws.Range("G3"))
$ws.Range("J2").Copy($ws.Range("J3"))
$ws.Range("L2").Copy($ws.Range("L3"))
$ws.Range("N2").Copy($ws.Range("N3"))
$ws.Range("P2").Copy($ws.Range("P3"))

$ws.Range("A3").Value = "PetPost_1"

$ws.Hyperlinks.Add($ws.Range("C3"), "https://live.virtualandemo.com/api/pets/findByTags?tags=grey")

# Restore C3's style to match C2 (undo the auto hyperlink style applied by Hyperlinks.Add)
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("B8").Select()
